# Update the title/date paragraph (first paragraph in the document)
$d = $word.ActiveDocument
$d.Paragraphs.Item(1).Range.Text = "2025-01-19 Sunday"

# Update the practice-problem table, one cell at a time (row/col addressed
# instead of a global Find/Replace, because several cells share identical
# "before" text but diverge in their replacement text).
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "20÷8=2, 4"   # was: 54÷6=9, 0
$t.Cell(1,2).Range.Text = "66÷4=16, 2"   # was: 26÷2=13, 0
$t.Cell(1,3).Range.Text = "10÷4=2, 2"   # was: 21÷9=2, 3
$t.Cell(1,4).Range.Text = "37÷7=5, 2"   # was: 31÷5=6, 1
$t.Cell(1,5).Range.Text = "21÷6=3, 3"   # was: 76÷7=10, 6
$t.Cell(5,1).Range.Text = "34÷6=5, 4"   # was: 28÷5=5, 3
$t.Cell(5,2).Range.Text = "87÷9=9, 6"   # was: 53÷9=5, 8
$t.Cell(5,3).Range.Text = "25÷3=8, 1"   # was: 52÷5=10, 2
$t.Cell(5,4).Range.Text = "87÷8=10, 7"   # was: 58÷3=19, 1
$t.Cell(5,5).Range.Text = "19÷2=9, 1"   # was: 83÷5=16, 3
$t.Cell(9,1).Range.Text = "53÷9=5, 8"   # was: 13÷7=1, 6
$t.Cell(9,2).Range.Text = "90÷2=45, 0"   # was: 67÷2=33, 1
$t.Cell(9,3).Range.Text = "15÷6=2, 3"   # was: 62÷8=7, 6
$t.Cell(9,4).Range.Text = "51÷4=12, 3"   # was: 16÷2=8, 0
$t.Cell(9,5).Range.Text = "72÷3=24, 0"   # was: 18÷8=2, 2
$t.Cell(13,1).Range.Text = "84÷4=21, 0"   # was: 79÷5=15, 4
$t.Cell(13,2).Range.Text = "49÷7=7, 0"   # was: 59÷4=14, 3
$t.Cell(13,3).Range.Text = "98÷4=24, 2"   # was: 32÷4=8, 0
$t.Cell(13,4).Range.Text = "85÷7=12, 1"   # was: 86÷3=28, 2
$t.Cell(13,5).Range.Text = "67÷4=16, 3"   # was: 57÷5=11, 2
$t.Cell(17,1).Range.Text = "97÷8=12, 1"   # was: 62÷8=7, 6
$t.Cell(17,2).Range.Text = "28÷3=9, 1"   # was: 39÷9=4, 3
$t.Cell(17,3).Range.Text = "12÷7=1, 5"   # was: 35÷9=3, 8
$t.Cell(17,4).Range.Text = "72÷8=9, 0"   # was: 22÷4=5, 2
$t.Cell(17,5).Range.Text = "46÷4=11, 2"   # was: 98÷9=10, 8
